$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, pushing the existing row 72 (and below) down to 73.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly record.
$ws.Cells.Item(72, 1).Value = 11
$ws.Cells.Item(72, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(72, 3).Value = "Bíobío"
$ws.Cells.Item(72, 4).Value = (Get-Date -Year 2022 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(72, 5).Value = 8
$ws.Cells.Item(72, 6).Value = 100112001
$ws.Cells.Item(72, 7).Value = "Berenjena"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 50
$ws.Cells.Item(72, 11).Value = 12000
$ws.Cells.Item(72, 12).Value = 13000
$ws.Cells.Item(72, 13).Value = 12400
$ws.Cells.Item(72, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(72, 15).Value = "Región Metropolitana"
$ws.Cells.Item(72, 16).Value = 207
$ws.Cells.Item(72, 17).Value = 60
$ws.Cells.Item(72, 18).Value = "Hortaliza"
